$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.062039498880902
$ws.Range("D2").Value = 1.060899797914149
$ws.Range("E2").Value = 1.065227820450102
$ws.Range("F2").Value = 1.070366925779918
$ws.Range("I2").Value = 1.027760580862774
$ws.Range("J2").Value = 1.067012076820579
$ws.Range("K2").Value = 1.063625321556726
$ws.Range("L2").Value = 1.06794162736605
$ws.Range("M2").Value = 1.073066953298517
$ws.Range("N2").Value = 1.06852735635766

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.065498731364662
$ws.Range("D3").Value = 1.064061909235623
$ws.Range("E3").Value = 1.068350223308369
$ws.Range("F3").Value = 1.073731499996781
$ws.Range("I3").Value = 1.028098158500023
$ws.Range("J3").Value = 1.070113658257546
$ws.Range("K3").Value = 1.066596039514865
$ws.Range("L3").Value = 1.070873631715965
$ws.Range("M3").Value = 1.076241584240708
$ws.Range("N3").Value = 1.071633342396024

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.067707034638981
$ws.Range("D4").Value = 1.066080400508695
$ws.Range("E4").Value = 1.070340970312035
$ws.Range("F4").Value = 1.075870741121192
$ws.Range("I4").Value = 1.028304821017455
$ws.Range("J4").Value = 1.072091078032288
$ws.Range("K4").Value = 1.06849060170308
$ws.Range("L4").Value = 1.072741083629925
$ws.Range("M4").Value = 1.078257889241557
$ws.Range("N4").Value = 1.07361357033366

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.06862841044887
$ws.Range("D5").Value = 1.066922546567713
$ws.Range("E5").Value = 1.071170965299074
$ws.Range("F5").Value = 1.07676122758596
$ws.Range("I5").Value = 1.02838891310831
$ws.Range("J5").Value = 1.072915505049633
$ws.Range("K5").Value = 1.069280621803646
$ws.Range("L5").Value = 1.073519216799699
$ws.Range("M5").Value = 1.079096679093159
$ws.Range("N5").Value = 1.074439168131932

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.068782708964313
$ws.Range("D6").Value = 1.067063574805807
$ws.Range("E6").Value = 1.071309924832661
$ws.Range("F6").Value = 1.07691023122941
$ws.Range("I6").Value = 1.028402869827211
$ws.Range("J6").Value = 1.073053531736199
$ws.Range("K6").Value = 1.069412896161508
$ws.Range("L6").Value = 1.073649466661693
$ws.Range("M6").Value = 1.079237001609822
$ws.Range("N6").Value = 1.074577390832225

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.067719373333542
$ws.Range("D7").Value = 1.066091678326313
$ws.Range("E7").Value = 1.070352087677126
$ws.Range("F7").Value = 1.075882674343561
$ws.Range("I7").Value = 1.028305955580028
$ws.Range("J7").Value = 1.07210212086123
$ws.Range("K7").Value = 1.068501183121605
$ws.Range("L7").Value = 1.072751508134857
$ws.Range("M7").Value = 1.078269131756515
$ws.Range("N7").Value = 1.073624628844685

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.063214907838218
$ws.Range("D8").Value = 1.061974275892833
$ws.Range("E8").Value = 1.066289305463967
$ws.Range("F8").Value = 1.071511971231653
$ws.Range("I8").Value = 1.027877122865667
$ws.Range("J8").Value = 1.068066494930347
$ws.Range("K8").Value = 1.064635131215362
$ws.Range("L8").Value = 1.068938782257584
$ws.Range("M8").Value = 1.074147808453698
$ws.Range("N8").Value = 1.069583271862087

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.055037920347559
$ws.Range("D9").Value = 1.054498924588384
$ws.Range("E9").Value = 1.058894327069477
$ws.Range("F9").Value = 1.063510135461595
$ws.Range("I9").Value = 1.027029869630035
$ws.Range("J9").Value = 1.060720490415538
$ws.Range("K9").Value = 1.057602314317832
$ws.Range("L9").Value = 1.061983985809324
$ws.Range("M9").Value = 1.066585508068448
$ws.Range("N9").Value = 1.062226835178269

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.049411762072623
$ws.Range("D10").Value = 1.049354940323541
$ws.Range("E10").Value = 1.053792983268282
$ws.Range("F10").Value = 1.057958910695703
$ws.Range("I10").Value = 1.026401297602511
$ws.Range("J10").Value = 1.055652561938521
$ws.Range("K10").Value = 1.052753492176318
$ws.Range("L10").Value = 1.057176191353614
$ws.Range("M10").Value = 1.061327844622511
$ws.Range("N10").Value = 1.057151709661515

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.046930836112355
$ws.Range("D11").Value = 1.047086516638081
$ws.Range("E11").Value = 1.051540360468183
$ws.Range("F11").Value = 1.055500227492829
$ws.Range("I11").Value = 1.026113475722859
$ws.Range("J11").Value = 1.053414594055205
$ws.Range("K11").Value = 1.050613001414127
$ws.Range("L11").Value = 1.055050793492065
$ws.Range("M11").Value = 1.058996514564633
$ws.Range("N11").Value = 1.054910563607106

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.046002289135736
$ws.Range("D12").Value = 1.046237487754847
$ws.Range("E12").Value = 1.050696794293894
$ws.Range("F12").Value = 1.054578382908434
$ws.Range("I12").Value = 1.026004166662352
$ws.Range("J12").Value = 1.052576499725328
$ws.Range("K12").Value = 1.049811519557837
$ws.Range("L12").Value = 1.05425451000126
$ws.Range("M12").Value = 1.05812202038228
$ws.Range("N12").Value = 1.05407127908715

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.04620178823427
$ws.Range("D13").Value = 1.046419903021902
$ws.Range("E13").Value = 1.050878056292221
$ws.Range("F13").Value = 1.05477651534342
$ws.Range("I13").Value = 1.026027723210624
$ws.Range("J13").Value = 1.052756586703609
$ws.Range("K13").Value = 1.04998373446759
$ws.Range("L13").Value = 1.054425628405317
$ws.Range("M13").Value = 1.058309993760818
$ws.Range("N13").Value = 1.054251621809597

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.046854227113049
$ws.Range("D14").Value = 1.047016468529913
$ws.Range("E14").Value = 1.051470772239406
$ws.Range("F14").Value = 1.055424204490692
$ws.Range("I14").Value = 1.026104489436773
$ws.Range("J14").Value = 1.053345457549602
$ws.Range("K14").Value = 1.050546882957989
$ws.Range("L14").Value = 1.054985113091988
$ws.Range("M14").Value = 1.058924404569219
$ws.Range("N14").Value = 1.054841328919735

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.047255276831293
$ws.Range("D15").Value = 1.047383171200181
$ws.Range("E15").Value = 1.051835049010283
$ws.Range("F15").Value = 1.055822120256954
$ws.Range("I15").Value = 1.026151468218869
$ws.Range("J15").Value = 1.053707368998242
$ws.Range("K15").Value = 1.050893000134044
$ws.Range("L15").Value = 1.055328918641205
$ws.Range("M15").Value = 1.059301822820091
$ws.Range("N15").Value = 1.055203754324144

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.049575444997431
$ws.Range("D16").Value = 1.049504600751038
$ws.Range("E16").Value = 1.05394153840813
$ws.Range("F16").Value = 1.058120900239524
$ws.Range("I16").Value = 1.026420065676154
$ws.Range("J16").Value = 1.055800148431029
$ws.Range("K16").Value = 1.052894665454154
$ws.Range("L16").Value = 1.057316306134181
$ws.Range("M16").Value = 1.061481388026169
$ws.Range("N16").Value = 1.057299505743771

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.051018637315916
$ws.Range("D17").Value = 1.050824143368992
$ws.Range("E17").Value = 1.055250992221271
$ws.Range("F17").Value = 1.059547926353391
$ws.Range("I17").Value = 1.026584327587963
$ws.Range("J17").Value = 1.05710105238935
$ws.Range("K17").Value = 1.054139122446916
$ws.Range("L17").Value = 1.058551087256864
$ws.Range("M17").Value = 1.062833704697344
$ws.Range("N17").Value = 1.058602257134929

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.051856128850742
$ws.Range("D18").Value = 1.05158986895532
$ws.Range("E18").Value = 1.056010576725281
$ws.Range("F18").Value = 1.060375005787006
$ws.Range("I18").Value = 1.026678632615597
$ws.Range("J18").Value = 1.057855667807477
$ws.Range("K18").Value = 1.054861063425544
$ws.Range("L18").Value = 1.059267127074994
$ws.Range("M18").Value = 1.063617229764502
$ws.Range("N18").Value = 1.059357944193505

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.052140971500413
$ws.Range("D19").Value = 1.051850301006893
$ws.Range("E19").Value = 1.056268871513451
$ws.Range("F19").Value = 1.060656132279898
$ws.Range("I19").Value = 1.026710534143189
$ws.Range("J19").Value = 1.058112271289886
$ws.Range("K19").Value = 1.055106567777926
$ws.Range("L19").Value = 1.059510575809933
$ws.Range("M19").Value = 1.063883508849801
$ws.Range("N19").Value = 1.059614912082292

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.050864243100124
$ws.Range("D20").Value = 1.050682978498503
$ws.Range("E20").Value = 1.055110936354836
$ws.Range("F20").Value = 1.059395368441283
$ws.Range("I20").Value = 1.026566859970232
$ws.Range("J20").Value = 1.056961911987049
$ws.Range("K20").Value = 1.054006012299225
$ws.Range("L20").Value = 1.058419041905478
$ws.Range("M20").Value = 1.062689160016047
$ws.Range("N20").Value = 1.058462919137297

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.046662296382254
$ws.Range("D21").Value = 1.046840974735528
$ws.Range("E21").Value = 1.051296423318575
$ws.Range("F21").Value = 1.055233715712502
$ws.Range("I21").Value = 1.026081950332214
$ws.Range("J21").Value = 1.053172240128662
$ws.Range("K21").Value = 1.050381228828609
$ws.Range("L21").Value = 1.054820549172362
$ws.Range("M21").Value = 1.058743714087457
$ws.Range("N21").Value = 1.054667865510189

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.043979588504024
$ws.Range("D22").Value = 1.044387977855773
$ws.Range("E22").Value = 1.048858355842908
$ws.Range("F22").Value = 1.052567304461876
$ws.Range("I22").Value = 1.025763160978784
$ws.Range("J22").Value = 1.050749955609426
$ws.Range("K22").Value = 1.048064967733674
$ws.Range("L22").Value = 1.052518455139343
$ws.Range("M22").Value = 1.056213510676665
$ws.Range("N22").Value = 1.052242141069096

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.045405710817983
$ws.Range("D23").Value = 1.045691994384902
$ws.Range("E23").Value = 1.050154683523342
$ws.Range("F23").Value = 1.053985653716602
$ws.Range("I23").Value = 1.025933492564635
$ws.Range("J23").Value = 1.052037900707913
$ws.Range("K23").Value = 1.049296479981049
$ws.Range("L23").Value = 1.053742683386482
$ws.Range("M23").Value = 1.057559624574062
$ws.Range("N23").Value = 1.053531915197352

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.050934020462288
$ws.Range("D24").Value = 1.050746776988141
$ws.Range("E24").Value = 1.055174234525287
$ws.Range("F24").Value = 1.0594643190914
$ws.Range("I24").Value = 1.026574757485799
$ws.Range("J24").Value = 1.057024796443368
$ws.Range("K24").Value = 1.054066171171111
$ws.Range("L24").Value = 1.058478720436655
$ws.Range("M24").Value = 1.062754489763549
$ws.Range("N24").Value = 1.058525892896756

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.057181629104521
$ws.Range("D25").Value = 1.056458798158188
$ws.Range("E25").Value = 1.060835329612098
$ws.Range("F25").Value = 1.065615848442761
$ws.Range("I25").Value = 1.027259971686622
$ws.Range("J25").Value = 1.062648698139201
$ws.Range("K25").Value = 1.059447787359118
$ws.Range("L25").Value = 1.063811204125454
$ws.Range("M25").Value = 1.068577531698473
$ws.Range("N25").Value = 1.064157781178068
